# Update invoice number, date, buyer details, item quantities/prices,
# and the "amount in words" strings (summ_str) as described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invoice number (merged G2:I2)
$ws.Range("G2").Value = 20

# Invoice date (merged A4:J4) - serial 43970 => 2020-05-19.
# Assign the raw date serial rather than a DateTime object so Excel does not
# reformat the cell with a new (time-aware) number format; the cell already
# carries a date-formatted style.
$ws.Range("A4").Value = 43970

# Buyer name (merged C9:J9)
$ws.Range("C9").Value = "ОАО ""пример"""

# Buyer extra detail (merged C10:J10)
$ws.Range("C10").Value = "бла-бла-бла"

# Line item 1: "Карта доступа Em-Marin белая с вырезом"
$ws.Range("D14").Value = 15
$ws.Range("F14").Value = 0.7

# Line item 2: "Пауч для ламинирования 80 мкм, 75х52 мм, уп. 100 шт."
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 23

# Amount in words, total incl. VAT (merged C21:J21)
$ws.Range("C21").Value = "сорок рублей ноль копеек"

# Amount in words, VAT only (merged C22:J22)
$ws.Range("C22").Value = "шесть рублей ноль копеек"
